$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header column labels (comment/label cleanup)
$ws.Range("A1").Value = "PubMed_title"
$ws.Range("B1").Value = "Formatted_title"

# Update the selected/visible range in the sheet view
$ws.Range("G10").Select()
